# Auto-generated Excel COM-interop script
# Applies updated price/profit figures to the Leve tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the "update Sheets via scheduled runner" commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2404.4
$ws.Range("I19").Value = 2374
$ws.Range("K19").Value = 2374
$ws.Range("M19").Value = -2199

$ws.Range("H33").Value = 151.66667
$ws.Range("I33").Value = 151.66667
$ws.Range("K33").Value = 151.66667
$ws.Range("M33").Value = 77.33332999999999

$ws.Range("H116").Value = 6000
$ws.Range("I116").Value = 6000
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = -2558
$ws.Range("N116").Value = -12884

$ws.Range("H132").Value = 4735.1665
$ws.Range("J132").Value = 1106
$ws.Range("L132").Value = 3318
$ws.Range("N132").Value = -8378

$ws.Range("H135").Value = 453.85715
$ws.Range("I135").Value = 196.16667
$ws.Range("K135").Value = 1765.50003
$ws.Range("M135").Value = 769.4999699999998

$ws.Range("H137").Value = 1759.8889
$ws.Range("I137").Value = 1724.8334
$ws.Range("J137").Value = 1830
$ws.Range("K137").Value = 5174.5002
$ws.Range("L137").Value = 5490
$ws.Range("M137").Value = -2624.5002
$ws.Range("N137").Value = -10590

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2532.8333
$ws.Range("I2").Value = 1439.4
$ws.Range("K2").Value = 1439.4
$ws.Range("M2").Value = -1326.4

$ws.Range("H88").Value = 2513.6875
$ws.Range("I88").Value = 927.5
$ws.Range("J88").Value = 4099.875
$ws.Range("K88").Value = 927.5
$ws.Range("L88").Value = 4099.875
$ws.Range("M88").Value = -521.5
$ws.Range("N88").Value = -4911.875

$ws.Range("H91").Value = 2513.6875
$ws.Range("I91").Value = 927.5
$ws.Range("J91").Value = 4099.875
$ws.Range("K91").Value = 927.5
$ws.Range("L91").Value = 4099.875
$ws.Range("M91").Value = 476.5
$ws.Range("N91").Value = -6907.875

$ws.Range("H97").Value = 1087.4
$ws.Range("I97").Value = 1087.4
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1087.4
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -591.4000000000001
$ws.Range("N97").Value = $null

$ws.Range("H116").Value = 2532.8333
$ws.Range("I116").Value = 1439.4
$ws.Range("K116").Value = 1439.4
$ws.Range("M116").Value = 854.5999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2532.8333
$ws.Range("I3").Value = 1439.4
$ws.Range("K3").Value = 1439.4
$ws.Range("M3").Value = -1325.4

$ws.Range("H22").Value = 392.30768
$ws.Range("I22").Value = 392.30768
$ws.Range("K22").Value = 392.30768
$ws.Range("M22").Value = -219.30768

$ws.Range("H64").Value = 908.1667
$ws.Range("I64").Value = 900.3333
$ws.Range("J64").Value = 912.0833
$ws.Range("K64").Value = 900.3333
$ws.Range("L64").Value = 912.0833
$ws.Range("M64").Value = -675.3333
$ws.Range("N64").Value = -1362.0833

$ws.Range("H67").Value = 908.1667
$ws.Range("I67").Value = 900.3333
$ws.Range("J67").Value = 912.0833
$ws.Range("K67").Value = 900.3333
$ws.Range("L67").Value = 912.0833
$ws.Range("M67").Value = -120.3333
$ws.Range("N67").Value = -2472.0833

$ws.Range("H94").Value = 2648.6667
$ws.Range("I94").Value = 2571.2727
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 2571.2727
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -2120.2727
$ws.Range("N94").Value = -4402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3805.724
$ws.Range("I31").Value = 3001.3
$ws.Range("J31").Value = 5593.3335
$ws.Range("K31").Value = 3001.3
$ws.Range("L31").Value = 5593.3335
$ws.Range("M31").Value = -2706.3
$ws.Range("N31").Value = -6183.3335

$ws.Range("H34").Value = 3805.724
$ws.Range("I34").Value = 3001.3
$ws.Range("J34").Value = 5593.3335
$ws.Range("K34").Value = 3001.3
$ws.Range("L34").Value = 5593.3335
$ws.Range("M34").Value = -2799.3
$ws.Range("N34").Value = -5997.3335

$ws.Range("H86").Value = 6000
$ws.Range("I86").Value = 6000
$ws.Range("K86").Value = 6000
$ws.Range("M86").Value = -4877

$ws.Range("H89").Value = 6000
$ws.Range("I89").Value = 6000
$ws.Range("K89").Value = 30000
$ws.Range("M89").Value = -24384

$ws.Range("H132").Value = 1386.9286
$ws.Range("I132").Value = 955.1539
$ws.Range("K132").Value = 2865.4617
$ws.Range("M132").Value = -335.4616999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 15000
$ws.Range("M57").Value = -14441

$ws.Range("H99").Value = 2229.8
$ws.Range("I99").Value = 2229.8
$ws.Range("K99").Value = 6689.400000000001
$ws.Range("M99").Value = -4443.400000000001

$ws.Range("H109").Value = 849.6
$ws.Range("J109").Value = 290
$ws.Range("L109").Value = 870
$ws.Range("N109").Value = -2950

$ws.Range("H129").Value = 1007.5
$ws.Range("I129").Value = 1007.5
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3022.5
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1977.5
$ws.Range("N129").Value = $null

$ws.Range("H132").Value = 1312.375
$ws.Range("J132").Value = 1312.375
$ws.Range("L132").Value = 11811.375
$ws.Range("N132").Value = -16871.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 822.94116
$ws.Range("I102").Value = 739.3333
$ws.Range("K102").Value = 739.3333
$ws.Range("M102").Value = 882.6667

$ws.Range("H122").Value = 11366287
$ws.Range("I122").Value = 13890572
$ws.Range("K122").Value = 41671716
$ws.Range("M122").Value = -41669266

$ws.Range("H132").Value = 3756
$ws.Range("I132").Value = 3756
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11268
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8738
$ws.Range("N132").Value = $null

$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -60119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2945.7778
$ws.Range("I122").Value = 3002
$ws.Range("J122").Value = 2900.8
$ws.Range("K122").Value = 9006
$ws.Range("L122").Value = 8702.400000000001
$ws.Range("M122").Value = -6556
$ws.Range("N122").Value = -13602.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1000
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

$ws.Range("H132").Value = 1955.3334
$ws.Range("I132").Value = 1742.5714
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 5227.7142
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -2697.7142
$ws.Range("N132").Value = -13160
